$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# D1 changes from "Cause of Error" to "Fix Implemented"; the old E1
# ("Fix Implemented") column is removed entirely.
$ws.Range("A1").Value = "Test Name"
$ws.Range("B1").Value = "Functionality Tested"
$ws.Range("C1").Value = "Failure Description"
$ws.Range("D1").Value = "Fix Implemented"
$ws.Range("E1").ClearContents()

# --- New data rows ---
# Column D is written before column C on each row so that shared-string
# indices get allocated in the same order as the target workbook
# (Login, Data validation, then the two failure-description strings).
$ws.Range("B2").Value = "Login"
$ws.Range("D2").Value = "Data validation"
$ws.Range("C2").Value = "No check to prevent multiple logins"

$ws.Range("B3").Value = "Login"
$ws.Range("D3").Value = "Data validation"
$ws.Range("C3").Value = "No restriction for admin-only transactions with standard login"

$ws.Range("B4").Value = "Login"

# --- Column widths ---
# Column 3 now holds the long failure-description text and column 4 the
# shorter fix-implemented text, so they are widened to fit.
$ws.Columns.Item(3).ColumnWidth = 48.333333333333336
$ws.Columns.Item(4).ColumnWidth = 27

# --- Selection ---
$ws.Range("C4").Select() | Out-Null
